$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Fill in new attendance marks ("p") for columns M and N -----------------
# Rows that get BOTH M and N filled in.
$rowsBoth = @(3,4,6,7,8,9,10,11,12,14,15,16,17,18,19,20,21,22,23,25,26,27,28,29,30,31,32,33)
foreach ($r in $rowsBoth) {
    $ws.Cells.Item($r, 13).Value = "p"   # column M
    $ws.Cells.Item($r, 14).Value = "p"   # column N
}

# Rows that only get column M filled in.
$rowsOnlyM = @(5,24)
foreach ($r in $rowsOnlyM) {
    $ws.Cells.Item($r, 13).Value = "p"   # column M
}

# Rows that only get column N filled in.
$rowsOnlyN = @(13,34)
foreach ($r in $rowsOnlyN) {
    $ws.Cells.Item($r, 14).Value = "p"   # column N
}

# --- Hide the now-unused columns --------------------------------------------
# Column D (the second "Alumnos/Correos" helper column) is hidden (width -> 0).
$ws.Columns("D").ColumnWidth = 0
$ws.Columns("D").Hidden = $true
# Columns F through J (attendance days that are no longer shown) are hidden (width -> 0).
$ws.Columns("F:J").ColumnWidth = 0
$ws.Columns("F:J").Hidden = $true

# --- Update the view / selection state --------------------------------------
$ws.Range("N34").Select()
